# Auto-update gym prices
#  - "4x4 Squat Racks": C2 price  $2,152.00 -> $2,149.00
#  - "Squat Stands":    C2 price  $1,556.00 -> $1,554.00
#                       and the now-stale "Rogue SM-2.5 Monster Squat Stand 2.0"
#                       row (row 6) is removed entirely.

$wb = $excel.ActiveWorkbook

# --- "4x4 Squat Racks" ------------------------------------------------------
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")
# Force text formatting first so the "$#,##0.00"-looking string is kept as a
# literal string instead of being parsed into a currency number.
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "$2,149.00"
# Drop back to the default cell style so no stray "@" number-format style
# lingers on the cell.
$ws1.Range("C2").Style = "Normal"

# --- "Squat Stands" ----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Squat Stands")
$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "$1,554.00"
$ws2.Range("C2").Style = "Normal"

# Remove the hyperlinks tied to the row-6 cells before deleting the row so no
# orphaned relationships are left behind, then delete the whole row, which
# shifts the dimension from A1:F6 down to A1:F5.
$ws2.Range("E6").Hyperlinks.Delete()
$ws2.Range("F6").Hyperlinks.Delete()
$ws2.Rows.Item(6).Delete()
